$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing rows 6..91 down to 7..92
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly record
$ws.Range('A6').Value = 2
$ws.Range('B6').Value = 'Comercializadora del Agro de Limarí'
$ws.Range('C6').Value = 'Coquimbo'
$ws.Range('D6').Value = 44699
$ws.Range('E6').Value = 4
$ws.Range('F6').Value = 100112030
$ws.Range('G6').Value = 'Poroto granado'
$ws.Range('H6').Value = 'Sin especificar'
$ws.Range('I6').Value = 'Primera'
$ws.Range('J6').Value = 360
$ws.Range('K6').Value = 16000
$ws.Range('L6').Value = 17000
$ws.Range('M6').Value = 16500
$ws.Range('N6').Value = '$/malla 25 kilos'
$ws.Range('O6').Value = 'Provincia de Limarí'
$ws.Range('P6').Value = 660
$ws.Range('Q6').Value = 25
$ws.Range('R6').Value = 'Hortaliza'
